$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# These D/E columns store plain text values (prices / % changes) as
# inline strings, not numbers. Force Text number format per-cell first so
# Excel does not auto-convert strings like "37.60" or "-0.60%" into numeric
# values/percentages on assignment.

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.60%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.08%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.161"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.02%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.929"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.75%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.23%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.001"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.21%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9316"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.78%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1094"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-11.21%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1922"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.76%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09053"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.51%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03305"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.43%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09596"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.00%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001393"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.88%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005746"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.23%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.595"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.76%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.434"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.25%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.999"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19.32%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.62%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2591"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.01%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04411"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.26%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001234"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.72%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004630"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "9.01%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001361"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.74%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003993"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-98.10%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02247"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.09%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05106"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.71%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007453"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-5.04%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009000"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-10.21%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1354"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.81%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002131"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.40%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008625"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-11.05%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006651"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.92%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002861"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-10.55%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001001"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-40.73%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.04%"
